# Auto-generated edit script: updates market-price derived columns (H-N)
# across multiple worksheets to reflect refreshed Universalis market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 604.7368
$ws.Range("J28").Value = 1062.8334
$ws.Range("L28").Value = 1062.8334
$ws.Range("N28").Value = -2032.8334
$ws.Range("H40").Value = 3837.7273
$ws.Range("I40").Value = 3439.5
$ws.Range("J40").Value = 4899.6665
$ws.Range("K40").Value = 3439.5
$ws.Range("L40").Value = 4899.6665
$ws.Range("M40").Value = -3264.5
$ws.Range("N40").Value = -5249.6665
$ws.Range("H74").Value = 4695.6
$ws.Range("I74").Value = 3478
$ws.Range("K74").Value = 3478
$ws.Range("M74").Value = -2542
$ws.Range("H77").Value = 4695.6
$ws.Range("I77").Value = 3478
$ws.Range("K77").Value = 17390
$ws.Range("M77").Value = -12710
$ws.Range("H86").Value = 1010131.6
$ws.Range("I86").Value = 2011983.4
$ws.Range("J86").Value = 8279.9
$ws.Range("K86").Value = 2011983.4
$ws.Range("L86").Value = 8279.9
$ws.Range("M86").Value = -2010860.4
$ws.Range("N86").Value = -10525.9
$ws.Range("H89").Value = 1010131.6
$ws.Range("I89").Value = 2011983.4
$ws.Range("J89").Value = 8279.9
$ws.Range("K89").Value = 10059917
$ws.Range("L89").Value = 41399.5
$ws.Range("M89").Value = -10054301
$ws.Range("N89").Value = -52631.5
$ws.Range("H106").Value = 9525.177
$ws.Range("I106").Value = 2592.8
$ws.Range("K106").Value = 2592.8
$ws.Range("M106").Value = -1961.8
$ws.Range("H113").Value = 125002500
$ws.Range("I113").Value = 50002000
$ws.Range("J113").Value = 200002990
$ws.Range("K113").Value = 50002000
$ws.Range("L113").Value = 200002990
$ws.Range("M113").Value = -49998746
$ws.Range("N113").Value = -200009498
$ws.Range("H138").Value = 3055.4792
$ws.Range("I138").Value = 1171.75
$ws.Range("J138").Value = 3432.225
$ws.Range("K138").Value = 3515.25
$ws.Range("L138").Value = 10296.675
$ws.Range("M138").Value = 1624.75
$ws.Range("N138").Value = -20576.675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1399.6897
$ws.Range("I2").Value = 914.2727
$ws.Range("K2").Value = 914.2727
$ws.Range("M2").Value = -801.2727
$ws.Range("H32").Value = 10419129
$ws.Range("I32").Value = 11906361
$ws.Range("J32").Value = 8499.833000000001
$ws.Range("K32").Value = 11906361
$ws.Range("L32").Value = 8499.833000000001
$ws.Range("M32").Value = -11906074
$ws.Range("N32").Value = -9073.833000000001
$ws.Range("H74").Value = 8071639.5
$ws.Range("I74").Value = 12501879
$ws.Range("K74").Value = 12501879
$ws.Range("M74").Value = -12501005
$ws.Range("H77").Value = 8071639.5
$ws.Range("I77").Value = 12501879
$ws.Range("K77").Value = 62509395
$ws.Range("M77").Value = -62505027
$ws.Range("H97").Value = 1638.9
$ws.Range("I97").Value = 1604.1875
$ws.Range("K97").Value = 1604.1875
$ws.Range("M97").Value = -1108.1875
$ws.Range("H116").Value = 1399.6897
$ws.Range("I116").Value = 914.2727
$ws.Range("K116").Value = 914.2727
$ws.Range("M116").Value = 1379.7273
$ws.Range("H122").Value = 2821.6155
$ws.Range("I122").Value = 1383
$ws.Range("K122").Value = 4149
$ws.Range("M122").Value = -1699

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1399.6897
$ws.Range("I3").Value = 914.2727
$ws.Range("K3").Value = 914.2727
$ws.Range("M3").Value = -800.2727
$ws.Range("H22").Value = 129.5
$ws.Range("I22").Value = 129.5
$ws.Range("K22").Value = 129.5
$ws.Range("M22").Value = 43.5
$ws.Range("H107").Value = 2334.7334
$ws.Range("I107").Value = 1646.6666
$ws.Range("K107").Value = 1646.6666
$ws.Range("M107").Value = 273.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 495.8
$ws.Range("I22").Value = 495.8
$ws.Range("K22").Value = 495.8
$ws.Range("M22").Value = -145.8
$ws.Range("H58").Value = 2150
$ws.Range("I58").Value = 1183.3334
$ws.Range("J58").Value = 2875
$ws.Range("K58").Value = 1183.3334
$ws.Range("L58").Value = 2875
$ws.Range("M58").Value = -980.3334
$ws.Range("N58").Value = -3281
$ws.Range("H110").Value = 99123.664
$ws.Range("J110").Value = 99123.664
$ws.Range("L110").Value = 99123.664
$ws.Range("N110").Value = -107303.664
$ws.Range("H122").Value = 2327.9333
$ws.Range("I122").Value = 2327.9333
$ws.Range("K122").Value = 6983.7999
$ws.Range("M122").Value = -4533.7999
$ws.Range("H132").Value = 1766.6957
$ws.Range("I132").Value = 1574.2727
$ws.Range("K132").Value = 4722.8181
$ws.Range("M132").Value = -2192.8181
$ws.Range("H134").Value = 358480.56
$ws.Range("I134").Value = 358480.56
$ws.Range("K134").Value = 1075441.68
$ws.Range("M134").Value = -1072906.68
$ws.Range("H136").Value = 2150
$ws.Range("I136").Value = 1183.3334
$ws.Range("J136").Value = 2875
$ws.Range("K136").Value = 3550.0002
$ws.Range("L136").Value = 8625
$ws.Range("M136").Value = -1000.0002
$ws.Range("N136").Value = -13725

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 252.46666
$ws.Range("I40").Value = 197.90909
$ws.Range("K40").Value = 791.63636
$ws.Range("M40").Value = -722.63636

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1139.8
$ws.Range("I2").Value = 174.75
$ws.Range("K2").Value = 174.75
$ws.Range("M2").Value = -61.75
$ws.Range("H70").Value = 5141
$ws.Range("I70").Value = 5045.5557
$ws.Range("K70").Value = 5045.5557
$ws.Range("M70").Value = -4775.5557
$ws.Range("H73").Value = 5141
$ws.Range("I73").Value = 5045.5557
$ws.Range("K73").Value = 5045.5557
$ws.Range("M73").Value = -4109.5557
$ws.Range("H122").Value = 1283.6154
$ws.Range("I122").Value = 988.625
$ws.Range("K122").Value = 2965.875
$ws.Range("M122").Value = -515.875
$ws.Range("H126").Value = 3875.7058
$ws.Range("I126").Value = 3220.889
$ws.Range("K126").Value = 9662.667000000001
$ws.Range("M126").Value = -7192.667000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1721.2142
$ws.Range("J16").Value = 2099.5
$ws.Range("L16").Value = 2099.5
$ws.Range("N16").Value = -2439.5
$ws.Range("H46").Value = 5931.5
$ws.Range("I46").Value = 3792.3333
$ws.Range("J46").Value = 8070.6665
$ws.Range("K46").Value = 3792.3333
$ws.Range("L46").Value = 8070.6665
$ws.Range("M46").Value = -3604.3333
$ws.Range("N46").Value = -8446.666499999999
$ws.Range("H82").Value = 825.6667
$ws.Range("I82").Value = 552.625
$ws.Range("K82").Value = 552.625
$ws.Range("M82").Value = -191.625
$ws.Range("H85").Value = 825.6667
$ws.Range("I85").Value = 552.625
$ws.Range("K85").Value = 552.625
$ws.Range("M85").Value = 695.375
$ws.Range("H93").Value = 66675670
$ws.Range("I93").Value = 66675670
$ws.Range("K93").Value = 66675670
$ws.Range("M93").Value = -66674422
$ws.Range("H140").Value = 99999.75
$ws.Range("J140").Value = 99999.75
$ws.Range("L140").Value = 99999.75
$ws.Range("N140").Value = -110359.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1459.25
$ws.Range("I100").Value = 1496.2858
$ws.Range("K100").Value = 2992.5716
$ws.Range("M100").Value = -2451.5716
$ws.Range("H122").Value = 15742.833
$ws.Range("I122").Value = 10927.25
$ws.Range("J122").Value = 25374
$ws.Range("K122").Value = 32781.75
$ws.Range("L122").Value = 76122
$ws.Range("M122").Value = -30331.75
$ws.Range("N122").Value = -81022
$ws.Range("H132").Value = 1600.8
$ws.Range("I132").Value = 1545.4286
$ws.Range("J132").Value = 1794.6
$ws.Range("K132").Value = 4636.2858
$ws.Range("L132").Value = 5383.799999999999
$ws.Range("M132").Value = -2106.2858
$ws.Range("N132").Value = -10443.8
$ws.Range("H136").Value = 697.7857
$ws.Range("I136").Value = 732.61536
$ws.Range("J136").Value = 245
$ws.Range("K136").Value = 2197.84608
$ws.Range("L136").Value = 735
$ws.Range("M136").Value = 352.1539199999997
$ws.Range("N136").Value = -5835
$ws.Range("H140").Value = 59428.5
$ws.Range("J140").Value = 59428.5
$ws.Range("L140").Value = 59428.5
$ws.Range("N140").Value = -69788.5
